$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 2-9 (columns A-I), reflecting the
# reordering + renaming of "Palestra" -> "Curso" + updated
# "Teoria da Comunicação" row data described by the commit.
$data = @(
    @("kjhlkjh",    "Baixa", 52, "Biblioteconomia",                     "Avançado",       2,  "Seminário", "brunocordeiro180", 100),
    @("C",          "Média", 9,  "Teoria da Comunicação",                "Intermediário",  17, "Workshop",  "brunocordeiro180", 956),
    @("Tópicos em Serviço Social de Educação", "Média", 88, "Serviço Social da Educação", "Básico", 5, "Workshop", "brunocordeiro180", 250),
    @("owo",         "Média", 2,  "Componentes da Dinâmica Demográfica", "Intermediário",  3,  "Curso",     "brunocordeiro180", 100),
    @("dfasadsf",    "Média", 21, "Contabilidade Nacional",               "Avançado",       2,  "Curso",     "brunocordeiro180", 456),
    @("nupcidade",   "Média", 2,  "Nupcialidade e Família",               "Básico",         2,  "Seminário", "brunocordeiro180", 567.23),
    @("kgkkjhkjh",   "Alta",  2,  "Conflitos e Coalizões Políticas",      "Avançado",       8,  "Curso",     "brunocordeiro180", 852),
    @("asdafd",      "Média", 8,  "Inflação",                             "Básico",         22, "Curso",     "brunocordeiro180", 100)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
}
